$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

# Rename the sheet (reflects "Through 2022-11-05" -> "Through 2022-11-07")
$ws.Name = "Through 2022-11-07"

# Update the row label text for November
$ws.Range("A12").Value = "November (through 11-07)"

# Update November row (row 12) values for columns B..I
$ws.Range("B12").Value = 8
$ws.Range("C12").Value = 17
$ws.Range("D12").Value = 24
$ws.Range("E12").Value = 18
$ws.Range("F12").Value = 10
$ws.Range("G12").Value = 45
$ws.Range("H12").Value = 50
$ws.Range("I12").Value = 18

# Update Total row (row 13) values for columns B..I
$ws.Range("B13").Value = 266
$ws.Range("C13").Value = 503
$ws.Range("D13").Value = 734
$ws.Range("E13").Value = 633
$ws.Range("F13").Value = 492
$ws.Range("G13").Value = 1102
$ws.Range("H13").Value = 1491
$ws.Range("I13").Value = 1418
